$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 157.64285
$ws.Range("I9").Value = 168.22223
$ws.Range("J9").Value = 138.6
$ws.Range("K9").Value = 168.22223
$ws.Range("L9").Value = 138.6
$ws.Range("M9").Value = 0.7777700000000038
$ws.Range("N9").Value = -476.6
$ws.Range("H19").Value = 815.1429000000001
$ws.Range("J19").Value = 651.375
$ws.Range("L19").Value = 651.375
$ws.Range("N19").Value = -1001.375
$ws.Range("H38").Value = 2826.8
$ws.Range("I38").Value = 59.666668
$ws.Range("J38").Value = 6977.5
$ws.Range("K38").Value = 179.000004
$ws.Range("L38").Value = 20932.5
$ws.Range("M38").Value = 192.999996
$ws.Range("N38").Value = -21676.5
$ws.Range("H58").Value = 1921.3334
$ws.Range("I58").Value = 1355.8889
$ws.Range("K58").Value = 4067.6667
$ws.Range("M58").Value = -3917.6667
$ws.Range("H62").Value = 8035.909
$ws.Range("I62").Value = 2347.25
$ws.Range("K62").Value = 2347.25
$ws.Range("M62").Value = -1723.25
$ws.Range("H65").Value = 8035.909
$ws.Range("I65").Value = 2347.25
$ws.Range("K65").Value = 11736.25
$ws.Range("M65").Value = -8616.25
$ws.Range("H80").Value = 883.3333
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("H83").Value = 883.3333
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("H127").Value = 619.6667
$ws.Range("I127").Value = 619.6667
$ws.Range("K127").Value = 1859.0001
$ws.Range("M127").Value = 3100.9999
$ws.Range("H131").Value = 2362.5
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("M83").ClearContents()
$ws.Range("N131").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 11865.315
$ws.Range("J44").Value = 11865.315
$ws.Range("L44").Value = 11865.315
$ws.Range("N44").Value = -12841.315
$ws.Range("H45").Value = 2135.25
$ws.Range("I45").Value = 1382.7142
$ws.Range("J45").Value = 3188.8
$ws.Range("K45").Value = 1382.7142
$ws.Range("L45").Value = 3188.8
$ws.Range("M45").Value = -1005.7142
$ws.Range("N45").Value = -3942.8
$ws.Range("H55").Value = 50333
$ws.Range("J55").Value = 99999
$ws.Range("L55").Value = 99999
$ws.Range("N55").Value = -100629
$ws.Range("H61").Value = 5843.846
$ws.Range("I61").Value = 4441.1113
$ws.Range("K61").Value = 4441.1113
$ws.Range("M61").Value = -4229.1113
$ws.Range("H97").Value = 613.9091
$ws.Range("I97").Value = 641.8
$ws.Range("K97").Value = 641.8
$ws.Range("M97").Value = -145.8
$ws.Range("H102").Value = 2792.5
$ws.Range("I102").Value = 1290.3572
$ws.Range("J102").Value = 8050
$ws.Range("K102").Value = 1290.3572
$ws.Range("L102").Value = 8050
$ws.Range("M102").Value = 331.6428000000001
$ws.Range("N102").Value = -11294
$ws.Range("H122").Value = 1205
$ws.Range("I122").Value = 1146
$ws.Range("K122").Value = 3438
$ws.Range("M122").Value = -988
$ws.Range("H136").Value = 5843.846
$ws.Range("I136").Value = 4441.1113
$ws.Range("K136").Value = 13323.3339
$ws.Range("M136").Value = -10773.3339

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 3725.8333
$ws.Range("I54").Value = 3271
$ws.Range("K54").Value = 3271
$ws.Range("M54").Value = -2787
$ws.Range("H86").Value = 3581.6365
$ws.Range("I86").Value = 1355.3846
$ws.Range("J86").Value = 6797.3335
$ws.Range("K86").Value = 1355.3846
$ws.Range("L86").Value = 6797.3335
$ws.Range("M86").Value = -232.3846000000001
$ws.Range("N86").Value = -9043.333500000001
$ws.Range("H89").Value = 3581.6365
$ws.Range("I89").Value = 1355.3846
$ws.Range("J89").Value = 6797.3335
$ws.Range("K89").Value = 6776.923000000001
$ws.Range("L89").Value = 33986.6675
$ws.Range("M89").Value = -1160.923000000001
$ws.Range("N89").Value = -45218.6675
$ws.Range("H94").Value = 277.33334
$ws.Range("I94").Value = 271
$ws.Range("K94").Value = 271
$ws.Range("M94").Value = 180
$ws.Range("H99").Value = 2746.3635
$ws.Range("I99").Value = 2245.1428
$ws.Range("K99").Value = 2245.1428
$ws.Range("M99").Value = -747.1428000000001
$ws.Range("H105").Value = 1668.2727
$ws.Range("I105").Value = 1435.4
$ws.Range("K105").Value = 1435.4
$ws.Range("M105").Value = 311.5999999999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 8000
$ws.Range("I136").Value = 4000
$ws.Range("K136").Value = 12000
$ws.Range("M136").Value = -6900
$ws.Range("H137").Value = 1875
$ws.Range("I137").Value = 583.3333
$ws.Range("J137").Value = 5750
$ws.Range("K137").Value = 1749.9999
$ws.Range("L137").Value = 17250
$ws.Range("M137").Value = 3350.0001
$ws.Range("N137").Value = -27450

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("H18").Value = 14550
$ws.Range("I18").Value = 14100
$ws.Range("J18").Value = 15000
$ws.Range("K18").Value = 14100
$ws.Range("L18").Value = 15000
$ws.Range("M18").Value = -13807
$ws.Range("N18").Value = -15586
$ws.Range("H42").Value = 98997
$ws.Range("J42").Value = 98997
$ws.Range("L42").Value = 98997
$ws.Range("N42").Value = -99967
$ws.Range("H80").Value = 2617.6
$ws.Range("I80").Value = 2362.6667
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 2362.6667
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -1364.6667
$ws.Range("N80").Value = -4996
$ws.Range("H83").Value = 2617.6
$ws.Range("I83").Value = 2362.6667
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 11813.3335
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -6821.333500000001
$ws.Range("N83").Value = -24984
$ws.Range("H115").Value = 98997
$ws.Range("J115").Value = 98997
$ws.Range("L115").Value = 98997
$ws.Range("N115").Value = -101347
$ws.Range("H122").Value = 3765.125
$ws.Range("I122").Value = 3590.2856
$ws.Range("J122").Value = 4989
$ws.Range("K122").Value = 10770.8568
$ws.Range("L122").Value = 14967
$ws.Range("M122").Value = -8320.856800000001
$ws.Range("N122").Value = -19867
$ws.Range("H126").Value = 1559.8
$ws.Range("I126").Value = 1500
$ws.Range("K126").Value = 4500
$ws.Range("M126").Value = -2030
$ws.Range("H132").Value = 3957.1667
$ws.Range("I132").Value = 3748.6
$ws.Range("K132").Value = 11245.8
$ws.Range("M132").Value = -8715.799999999999
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4340.25
$ws.Range("I40").Value = 4453.6665
$ws.Range("K40").Value = 4453.6665
$ws.Range("M40").Value = -4317.6665
$ws.Range("H93").Value = 795.3333
$ws.Range("I93").Value = 795.3333
$ws.Range("K93").Value = 795.3333
$ws.Range("M93").Value = 452.6667
$ws.Range("H122").Value = 3003
$ws.Range("I122").Value = 3000
$ws.Range("J122").Value = 3003.4285
$ws.Range("K122").Value = 9000
$ws.Range("L122").Value = 9010.2855
$ws.Range("M122").Value = -6550
$ws.Range("N122").Value = -13910.2855
$ws.Range("H132").Value = 4941.5
$ws.Range("I132").Value = 4941.5
$ws.Range("K132").Value = 14824.5
$ws.Range("M132").Value = -12294.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7675.6924
$ws.Range("J62").Value = 8333.143
$ws.Range("L62").Value = 8333.143
$ws.Range("N62").Value = -9581.143
$ws.Range("H65").Value = 7675.6924
$ws.Range("J65").Value = 8333.143
$ws.Range("L65").Value = 41665.715
$ws.Range("N65").Value = -47905.715
$ws.Range("H96").Value = 1798.3334
$ws.Range("I96").Value = 1697.5
$ws.Range("K96").Value = 1697.5
$ws.Range("M96").Value = -324.5
$ws.Range("H100").Value = 925
$ws.Range("J100").Value = 1376
$ws.Range("L100").Value = 2752
$ws.Range("N100").Value = -3834
$ws.Range("H113").Value = 958.125
$ws.Range("I113").Value = 861
$ws.Range("K113").Value = 2583
$ws.Range("M113").Value = -413
$ws.Range("H122").Value = 3629.7273
$ws.Range("I122").Value = 2240.25
$ws.Range("J122").Value = 7335
$ws.Range("K122").Value = 6720.75
$ws.Range("L122").Value = 22005
$ws.Range("M122").Value = -4270.75
$ws.Range("N122").Value = -26905
$ws.Range("H132").Value = 2101.6667
$ws.Range("J132").Value = 2101.6667
$ws.Range("L132").Value = 6305.000100000001
$ws.Range("N132").Value = -11365.0001
